$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New adj_close (column G) text values for data rows 2..253, in row order.
# These come from the refreshed R-script pull; values are stored as text
# (shared strings) in the original workbook, so we must force text entry
# via a leading apostrophe and then strip the resulting quote-prefix style
# so the cell keeps the default (unstyled) look.
$gValues = @(
    "2.88316249847412","2.82511234283447","2.76706194877625","2.78641223907471","2.80576229095459","2.81543731689453","2.82511234283447","2.79608726501465","2.82511234283447","2.78834700584412",
    "2.8038272857666","2.79608726501465","2.80189204216003","2.807697057724","2.82511234283447","2.90057754516602","3.05731344223022","3.0282883644104","3.04763841629028","3.03796339035034",
    "2.999263048172","3.02441811561584","2.9818480014801","3.03989815711975","3.04957318305969","3.09988355636597","3.06311821937561","3.02054810523987","3.07472825050354","2.99539303779602",
    "2.95088768005371","2.93153762817383","2.91799259185791","2.89283752441406","2.86574745178223","2.86574745178223","2.82511234283447","2.8541374206543","2.83865737915039","2.8038272857666",
    "2.82124209403992","2.81156706809998","2.80576229095459","2.8038272857666","2.84446239471436","2.81930732727051","2.81350231170654","2.8541374206543","2.83478736877441","2.79802203178406",
    "2.77673697471619","2.72062182426453","2.81543731689453","2.78641223907471","2.79028224945068","2.78641223907471","2.80576229095459","2.84059238433838","2.80576229095459","2.8541374206543",
    "2.92573261260986","2.90057754516602","2.87348747253418","2.79221701622009","2.80576229095459","2.75351691246033","2.71288180351257","2.76512694358826","2.72836184501648","2.7244918346405",
    "2.67998671531677","2.77673697471619","2.76319193840027","2.72255682945251","2.77673697471619","2.77673697471619","2.74384188652039","2.72062182426453","2.71675181388855","2.71094679832458",
    "2.67998671531677","2.7090117931366","2.67031168937683","2.6451563835144","2.67418146133423","2.64709162712097","2.65676665306091","2.65483140945435","2.65096163749695","2.68579173088074",
    "2.67224669456482","2.7573869228363","2.75351691246033","2.7573869228363","2.71868681907654","2.75158214569092","2.7573869228363","2.73803687095642","2.81543731689453","2.84252738952637",
    "2.97604298591614","2.93153762817383","2.94508266448975","2.9489529132843","2.80576229095459","2.85994243621826","2.80963230133057","2.89283752441406","2.87929248809814","2.87155246734619",
    "2.902512550354","2.81543731689453","2.85994243621826","2.87155246734619","2.8541374206543","2.88896751403809","2.89090251922607","2.87929248809814","2.85994243621826","2.88316249847412",
    "2.90444779396057","2.85800743103027","2.8696174621582","2.8541374206543","2.89477276802063","2.92186260223389","2.88316249847412","2.85026741027832","2.83285236358643","2.81543731689453",
    "2.78641223907471","2.79995727539062","2.73803687095642","2.73029685020447","2.71481680870056","2.76319193840027","2.74771189689636","2.73997187614441","2.69933676719666","2.76512694358826",
    "2.77286696434021","2.75158214569092","2.72836184501648","2.7573869228363","2.75351691246033","2.77673697471619","2.7748019695282","2.80576229095459","2.81350231170654","2.83478736877441",
    "2.8541374206543","2.84639739990234","2.84446239471436","2.83091735839844","2.79415225982666","2.7244918346405","2.92960286140442","2.91025257110596","2.94121265411377","3.02635312080383",
    "3.01087307929993","2.97797799110413","2.97410774230957","2.99539303779602","3.08633852005005","3.13471364974976","3.12116861343384","3.16373872756958","3.21211385726929","3.24113893508911",
    "3.38819980621338","3.35723972320557","3.3378894329071","3.39206981658936","3.34175944328308","3.2895143032074","3.21017909049988","3.18502378463745","3.14632368087769","3.03022313117981",
    "2.98958802223206","3.03022313117981","3.05537819862366","3.10568857192993","3.09214353561401","3.03022313117981","3.04376816749573","3.01861310005188","3.06311821937561","3.07859826087952",
    "3.0650532245636","3.04763841629028","3.04570317268372","3.06118321418762","3.01377558708191","3.06215071678162","3.04763841629028","3.01861310005188","3.00410056114197","2.98475050926208",
    "2.97991275787354","2.97023773193359","2.95088768005371","2.87348747253418","2.85897493362427","2.83478736877441","2.86381244659424","2.83478736877441","2.85897493362427","2.76706194877625",
    "2.73803687095642","2.7090117931366","2.74771189689636","2.76706194877625","2.77673697471619","2.73803687095642","2.78157448768616","2.80576229095459","2.86381244659424","2.88800001144409",
    "2.90735030174255","2.88800001144409","2.9170253276825","2.93637537956238","2.94605040550232","3.01861310005188","2.98475050926208","3.01377558708191","2.999263048172","2.98475050926208",
    "2.96056270599365","2.97023773193359","2.999263048172","2.99442553520203","2.96056270599365","2.98475050926208","2.98958802223206","2.98475050926208","2.91218757629395","2.93637537956238",
    "2.96056270599365","2.98958802223206","3.00893807411194","2.97991275787354","2.99442553520203","2.92186260223389","2.85897493362427","2.86864995956421","2.8541374206543","2.86381244659424",
    "2.8541374206543","2.88800001144409"
)

for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 7)
    $cell.Value = "'" + $gValues[$i]
    $cell.Style = "Normal"
}

# Row 253 was the latest (partial-day) quote; it is now a completed trading
# day, so its timestamp collapses to the standard midnight fraction used by
# every other row.
$ws.Range("A253").Value = 45464.2916666667

# Append the new row 254 pulled by the refreshed R script, copying the date
# cell's number format from row 253 so it keeps the datetime display.
$ws.Range("A253:H253").Copy()
$ws.Range("A254:H254").PasteSpecial(-4122)

$ws.Cells.Item(254, 1).Value = 45467.6494791667
$ws.Cells.Item(254, 2).Value = 303766
$ws.Cells.Item(254, 3).Value = 2.96000003814697
$ws.Cells.Item(254, 4).Value = 2.89000010490417
$ws.Cells.Item(254, 5).Value = 2.97000002861023
$ws.Cells.Item(254, 6).Value = 2.88499999046326

$gNew = $ws.Cells.Item(254, 7)
$gNew.Value = "'2.88499999046326"
$gNew.Style = "Normal"

$ws.Cells.Item(254, 8).Value = "YACHT.MI"

